# Re-ran resolve and classify+summarise steps after changes to mapping file.
# Updates the summary tables on "Range Status", "Species qualification" and
# "High Priority break-up" sheets to reflect the freshly computed values.

$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet --------------------------------------------------
# Species (no.) column collapses to 0 for every range class, and the
# Species (perc.) column is no longer populated (cells removed).
$wsRange = $wb.Worksheets.Item("Range Status")

$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()

$wsRange.Range("C3").ClearContents()

$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()

$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()

$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()

$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- "Species qualification" sheet -----------------------------------------
# "Range Analysis" row now shows 0 species selected for analysis.
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- "High Priority break-up" sheet -----------------------------------------
# Updated high-priority species break-up percentages/counts.
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("E2").Value = 9.1
$wsBreakup.Range("D3").Value = 10
$wsBreakup.Range("E3").Value = 90.90000000000001
